# Commit message: "disable/ enable dropdown ribbon, explorer"
# This edit changes the Image Name for the Import dropdown (row 3) and the
# Browse project button (row 16) from their old icon names ("import_file"
# and "browse_project") to a shared new icon name "project_new_16x16".
# It also moves the active cell selection to E17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("ExplorerTab")

# Row 3 (ImportProjectDropdown): Image Name column E
$ws.Range("E3").Value = "project_new_16x16"

# Row 16 (BrowseProjectButton): Image Name column E
$ws.Range("E16").Value = "project_new_16x16"

# Update the selected/active cell as recorded in the saved view state
$ws.Range("E17").Select()
